$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.969.84"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "'2.407.00"
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'554.23"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "'142.17"
$ws.Range("E6").Value = "  +3.63%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").Value = "'2.400.74"
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("E14").Value = "  +4.12%  "
$ws.Range("E15").Value = "  +5.31%  "
$ws.Range("D16").Value = "'2.841.11"
$ws.Range("E16").Value = "  +2.36%  "
$ws.Range("D17").Value = "'61.968.87"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").Value = "'2.399.66"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").Value = "'11.07"
$ws.Range("E19").Value = "  +3.26%  "
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "'322.74"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").Value = "'6.70"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'65.03"
$ws.Range("E24").Value = "  +1.73%  "
$ws.Range("E25").Value = "  +4.03%  "
$ws.Range("D26").Value = "'9.01"
$ws.Range("E26").Value = "  +9.55%  "
$ws.Range("D27").Value = "'575.52"
$ws.Range("E27").Value = "  +16.15%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "'2.524.33"
$ws.Range("E29").Value = "  +2.01%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0₃0925"
$ws.Range("E30").Value = "  +5.73%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.23"
$ws.Range("E31").Value = "  +2.17%  "
$ws.Range("E32").Value = "  +5.64%  "
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("E34").Value = "  +2.64%  "
$ws.Range("E35").Value = "  +3.82%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "'5.63"
$ws.Range("E37").Value = "  +6.11%  "
$ws.Range("D38").Value = "'4.75"
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D40").Value = "'150.61"
$ws.Range("E40").Value = "  +3.70%  "
$ws.Range("D41").Value = "'18.59"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("E42").Value = "  -1.92%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  +13.63%  "
$ws.Range("D45").Value = "'149.35"
$ws.Range("E45").Value = "  +1.93%  "
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D48").Value = "'20.08"
$ws.Range("E48").Value = "  +5.22%  "
$ws.Range("D49").Value = "'0.587"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("E51").Value = "  +2.53%  "
